$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights -------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 43.5
$ws.Rows.Item(2).RowHeight = 101.5
$ws.Rows.Item(4).RowHeight = 72.5
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).RowHeight = 87

# --- Style: apply wrap-text alignment to the cells that need it --------
$ws.Range("A1").WrapText = $true
$ws.Range("C1").WrapText = $true
$ws.Range("D1").WrapText = $true
$ws.Range("C2").WrapText = $true
$ws.Range("A5").WrapText = $true
$ws.Range("E5").WrapText = $true
$ws.Range("G5").WrapText = $true
$ws.Range("B6").WrapText = $true

# --- Cell value updates --------------------------------------------------
# NOTE: string cells are written in the exact order the new shared-string
# table expects so the resulting uniqueCount/order matches (B6, G1, C2,
# F2, G3, C4, E4).

$ws.Range("B6").Value = "31`n세계 금연의 날`n바다의 날"

$ws.Range("F1").Value = 3
$ws.Range("G1").Value = "4`n지식재산의 날"

$ws.Range("A2").Value = 5
$ws.Range("C2").Value = "7`n음8.1`n백로`n푸른하늘의날"
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = "10`n세계 자살예방의 날`n9.10 해양경찰의날"
$ws.Range("G2").Value = 11

$ws.Range("B3").Value = 13
$ws.Range("F3").Value = 17
$ws.Range("G3").Value = "18`n청년의날"

$ws.Range("A4").Value = 19
$ws.Range("C4").Value = "21`n음8.15`n추석`n치매극복의 날"
$ws.Range("D4").Value = 22
$ws.Range("E4").Value = "23`n추분"

$ws.Range("B5").Value = 27

$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 5
